$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API")
$ws.Range("A1").Value = "Test"
